$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 3531.625
$ws.Cells.Item(100, 9).Value = 2786.25
$ws.Cells.Item(100, 10).Value = 4277
$ws.Cells.Item(100, 11).Value = 2786.25
$ws.Cells.Item(100, 12).Value = 4277
$ws.Cells.Item(100, 13).Value = -2245.25
$ws.Cells.Item(100, 14).Value = -5359
$ws.Cells.Item(113, 8).Value = 4415.75
$ws.Cells.Item(113, 9).Value = 4946.8335
$ws.Cells.Item(113, 10).Value = 2822.5
$ws.Cells.Item(113, 11).Value = 4946.8335
$ws.Cells.Item(113, 12).Value = 2822.5
$ws.Cells.Item(113, 13).Value = -1692.8335
$ws.Cells.Item(113, 14).Value = -9330.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(23, 8).Value = 26250
$ws.Cells.Item(23, 10).Value = 35000
$ws.Cells.Item(23, 12).Value = 35000
$ws.Cells.Item(23, 14).Value = -35518
$ws.Cells.Item(25, 8).Value = 10499.333
$ws.Cells.Item(25, 10).Value = 12999.667
$ws.Cells.Item(25, 12).Value = 12999.667
$ws.Cells.Item(25, 14).Value = -13803.667
$ws.Cells.Item(30, 8).Value = 3535.3333
$ws.Cells.Item(30, 10).Value = 3535.3333
$ws.Cells.Item(30, 12).Value = 3535.3333
$ws.Cells.Item(30, 14).Value = -3835.3333
$ws.Cells.Item(32, 8).Value = 1415.5977
$ws.Cells.Item(32, 9).Value = 1026.3494
$ws.Cells.Item(32, 11).Value = 1026.3494
$ws.Cells.Item(32, 13).Value = -739.3494000000001
$ws.Cells.Item(35, 8).Value = 5499.5
$ws.Cells.Item(35, 9).Value = 5499.5
$ws.Cells.Item(35, 11).Value = 5499.5
$ws.Cells.Item(35, 13).Value = -5093.5
$ws.Cells.Item(36, 8).Value = 2708.6667
$ws.Cells.Item(36, 9).Value = 2708.6667
$ws.Cells.Item(36, 11).Value = 2708.6667
$ws.Cells.Item(36, 13).Value = -2362.6667
$ws.Cells.Item(37, 8).Value = 22746.875
$ws.Cells.Item(37, 10).Value = 19737.5
$ws.Cells.Item(37, 12).Value = 19737.5
$ws.Cells.Item(37, 14).Value = -20283.5
$ws.Cells.Item(38, 8).Value = 10000
$ws.Cells.Item(38, 10).Value = 10000
$ws.Cells.Item(38, 12).Value = 10000
$ws.Cells.Item(38, 14).Value = -10934
$ws.Cells.Item(42, 8).Value = 3000
$ws.Cells.Item(42, 10).Value = 3000
$ws.Cells.Item(42, 12).Value = 3000
$ws.Cells.Item(42, 14).Value = -3972
$ws.Cells.Item(44, 8).Value = 76369.664
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 76369.664
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 76369.664
$ws.Cells.Item(44, 13).ClearContents()
$ws.Cells.Item(44, 14).Value = -77345.664
$ws.Cells.Item(122, 8).Value = 6815.3
$ws.Cells.Item(122, 9).Value = 6269.25
$ws.Cells.Item(122, 11).Value = 18807.75
$ws.Cells.Item(122, 13).Value = -16357.75
$ws.Cells.Item(132, 8).Value = 4730.426
$ws.Cells.Item(132, 9).Value = 5521.758
$ws.Cells.Item(132, 11).Value = 16565.274
$ws.Cells.Item(132, 13).Value = -14035.274

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3521.9
$ws.Cells.Item(99, 9).Value = 4124.875
$ws.Cells.Item(99, 11).Value = 4124.875
$ws.Cells.Item(99, 13).Value = -2626.875
$ws.Cells.Item(107, 8).Value = 3358.8
$ws.Cells.Item(107, 9).Value = 2933
$ws.Cells.Item(107, 11).Value = 2933
$ws.Cells.Item(107, 13).Value = -1013

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 58499.75
$ws.Cells.Item(68, 10).Value = 58499.75
$ws.Cells.Item(68, 12).Value = 58499.75
$ws.Cells.Item(68, 14).Value = -59997.75
$ws.Cells.Item(71, 8).Value = 58499.75
$ws.Cells.Item(71, 10).Value = 58499.75
$ws.Cells.Item(71, 12).Value = 175499.25
$ws.Cells.Item(71, 14).Value = -182987.25
$ws.Cells.Item(74, 8).Value = 36492.2
$ws.Cells.Item(74, 10).Value = 39597.75
$ws.Cells.Item(74, 12).Value = 39597.75
$ws.Cells.Item(74, 14).Value = -41345.75
$ws.Cells.Item(77, 8).Value = 36492.2
$ws.Cells.Item(77, 10).Value = 39597.75
$ws.Cells.Item(77, 12).Value = 118793.25
$ws.Cells.Item(77, 14).Value = -127529.25
$ws.Cells.Item(93, 8).Value = 13692.833
$ws.Cells.Item(93, 9).Value = 12631.4
$ws.Cells.Item(93, 11).Value = 12631.4
$ws.Cells.Item(93, 13).Value = -10759.4
$ws.Cells.Item(99, 8).Value = 1434
$ws.Cells.Item(99, 9).Value = 1438.8572
$ws.Cells.Item(99, 10).Value = 1400
$ws.Cells.Item(99, 11).Value = 1438.8572
$ws.Cells.Item(99, 12).Value = 1400
$ws.Cells.Item(99, 13).Value = 59.14280000000008
$ws.Cells.Item(99, 14).Value = -4396
$ws.Cells.Item(105, 8).Value = 1481.7
$ws.Cells.Item(105, 9).Value = 1550.3334
$ws.Cells.Item(105, 11).Value = 1550.3334
$ws.Cells.Item(105, 13).Value = 196.6666
$ws.Cells.Item(107, 8).Value = 1075.8334
$ws.Cells.Item(107, 9).Value = 1002.6667
$ws.Cells.Item(107, 10).Value = 1149
$ws.Cells.Item(107, 11).Value = 1002.6667
$ws.Cells.Item(107, 12).Value = 1149
$ws.Cells.Item(107, 13).Value = 917.3333
$ws.Cells.Item(107, 14).Value = -4989
$ws.Cells.Item(122, 8).Value = 1789.125
$ws.Cells.Item(122, 9).Value = 1883.5
$ws.Cells.Item(122, 11).Value = 5650.5
$ws.Cells.Item(122, 13).Value = -3200.5
$ws.Cells.Item(126, 8).Value = 1434
$ws.Cells.Item(126, 9).Value = 1438.8572
$ws.Cells.Item(126, 10).Value = 1400
$ws.Cells.Item(126, 11).Value = 4316.571599999999
$ws.Cells.Item(126, 12).Value = 4200
$ws.Cells.Item(126, 13).Value = -1846.571599999999
$ws.Cells.Item(126, 14).Value = -9140

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 540597.75
$ws.Cells.Item(2, 9).Value = 769260.75
$ws.Cells.Item(2, 11).Value = 4615564.5
$ws.Cells.Item(2, 13).Value = -4615451.5
$ws.Cells.Item(8, 8).Value = 982.3333
$ws.Cells.Item(8, 9).Value = 982.3333
$ws.Cells.Item(8, 11).Value = 2946.9999
$ws.Cells.Item(8, 13).Value = -2807.9999
$ws.Cells.Item(34, 8).Value = 1899
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(39, 8).Value = 7665.6665
$ws.Cells.Item(39, 10).Value = 7665.6665
$ws.Cells.Item(39, 12).Value = 22996.9995
$ws.Cells.Item(39, 14).Value = -23584.9995
$ws.Cells.Item(51, 8).Value = 1059.4
$ws.Cells.Item(51, 9).Value = 433
$ws.Cells.Item(51, 11).Value = 1299
$ws.Cells.Item(51, 13).Value = -839
$ws.Cells.Item(55, 8).Value = 6651.1333
$ws.Cells.Item(55, 9).Value = 846.375
$ws.Cells.Item(55, 10).Value = 13285.143
$ws.Cells.Item(55, 11).Value = 2539.125
$ws.Cells.Item(55, 12).Value = 39855.429
$ws.Cells.Item(55, 13).Value = -2362.125
$ws.Cells.Item(55, 14).Value = -40209.429
$ws.Cells.Item(57, 8).Value = 17999.5
$ws.Cells.Item(57, 9).Value = 17999.5
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 53998.5
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = -53439.5
$ws.Cells.Item(57, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 1579.0667
$ws.Cells.Item(68, 10).Value = 2272.75
$ws.Cells.Item(68, 12).Value = 6818.25
$ws.Cells.Item(68, 14).Value = -8440.25
$ws.Cells.Item(71, 8).Value = 1579.0667
$ws.Cells.Item(71, 10).Value = 2272.75
$ws.Cells.Item(71, 12).Value = 20454.75
$ws.Cells.Item(71, 14).Value = -28566.75
$ws.Cells.Item(129, 8).Value = 3849638.8
$ws.Cells.Item(129, 9).Value = 4512.25
$ws.Cells.Item(129, 10).Value = 7145461.5
$ws.Cells.Item(129, 11).Value = 13536.75
$ws.Cells.Item(129, 12).Value = 21436384.5
$ws.Cells.Item(129, 13).Value = -8536.75
$ws.Cells.Item(129, 14).Value = -21446384.5

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 23000
$ws.Cells.Item(26, 9).Value = 26000
$ws.Cells.Item(26, 10).Value = 20000
$ws.Cells.Item(26, 11).Value = 26000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = -25720
$ws.Cells.Item(26, 14).Value = -20560
$ws.Cells.Item(50, 8).Value = 23000
$ws.Cells.Item(50, 9).Value = 26000
$ws.Cells.Item(50, 10).Value = 20000
$ws.Cells.Item(50, 11).Value = 26000
$ws.Cells.Item(50, 12).Value = 20000
$ws.Cells.Item(50, 13).Value = -25502
$ws.Cells.Item(50, 14).Value = -20996
$ws.Cells.Item(57, 8).Value = 35599.6
$ws.Cells.Item(57, 9).Value = 32999
$ws.Cells.Item(57, 10).Value = 37333.332
$ws.Cells.Item(57, 11).Value = 32999
$ws.Cells.Item(57, 12).Value = 37333.332
$ws.Cells.Item(57, 13).Value = -32179
$ws.Cells.Item(57, 14).Value = -38973.332
$ws.Cells.Item(102, 8).Value = 3075.889
$ws.Cells.Item(102, 9).Value = 2886.7778
$ws.Cells.Item(102, 11).Value = 2886.7778
$ws.Cells.Item(102, 13).Value = -1264.7778
$ws.Cells.Item(107, 8).Value = 2161.1667
$ws.Cells.Item(107, 9).Value = 3024.5
$ws.Cells.Item(107, 10).Value = 1297.8334
$ws.Cells.Item(107, 11).Value = 3024.5
$ws.Cells.Item(107, 12).Value = 1297.8334
$ws.Cells.Item(107, 13).Value = -1104.5
$ws.Cells.Item(107, 14).Value = -5137.8334
$ws.Cells.Item(122, 8).Value = 100000
$ws.Cells.Item(122, 9).Value = 100000
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 300000
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -297550
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 12194068
$ws.Cells.Item(126, 9).Value = 7290.5713
$ws.Cells.Item(126, 10).Value = 21672672
$ws.Cells.Item(126, 11).Value = 21871.7139
$ws.Cells.Item(126, 12).Value = 65018016
$ws.Cells.Item(126, 13).Value = -19401.7139
$ws.Cells.Item(126, 14).Value = -65022956
$ws.Cells.Item(132, 8).Value = 8897.729499999999
$ws.Cells.Item(132, 9).Value = 8711.944
$ws.Cells.Item(132, 10).Value = 9455.083000000001
$ws.Cells.Item(132, 11).Value = 26135.832
$ws.Cells.Item(132, 12).Value = 28365.249
$ws.Cells.Item(132, 13).Value = -23605.832
$ws.Cells.Item(132, 14).Value = -33425.249

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3547.8823
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 3547.8823
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 3547.8823
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).Value = -4137.8823
$ws.Cells.Item(27, 8).Value = 3547.8823
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 3547.8823
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 3547.8823
$ws.Cells.Item(27, 13).ClearContents()
$ws.Cells.Item(27, 14).Value = -3761.8823
$ws.Cells.Item(122, 8).Value = 4790.1763
$ws.Cells.Item(122, 9).Value = 4826.4443
$ws.Cells.Item(122, 11).Value = 14479.3329
$ws.Cells.Item(122, 13).Value = -12029.3329

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(50, 8).Value = 36250
$ws.Cells.Item(50, 10).Value = 36250
$ws.Cells.Item(50, 12).Value = 36250
$ws.Cells.Item(50, 14).Value = -37512
$ws.Cells.Item(122, 8).Value = 13891982
$ws.Cells.Item(122, 10).Value = 2397
$ws.Cells.Item(122, 12).Value = 7191
$ws.Cells.Item(122, 14).Value = -12091
$ws.Cells.Item(126, 8).Value = 70179704
$ws.Cells.Item(126, 9).Value = 23812674
$ws.Cells.Item(126, 10).Value = 200007400
$ws.Cells.Item(126, 11).Value = 71438022
$ws.Cells.Item(126, 12).Value = 600022200
$ws.Cells.Item(126, 13).Value = -71435552
$ws.Cells.Item(126, 14).Value = -600027140
